$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109864/1750183612496/A5MTrC4t3W-camera-video-segment-driver-1750183614996.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNED233MW%2F20250618%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250618T160442Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIBxOo0S7hIeADNS8X%2F11v%2FFXQ5yY40M9wwoRF9Veb0znAiA0hEdIMmymyjsOlVstbTPCrvZ04Qw671YfXkdetNsKSCrmAwiR%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM14o7WDIZqelIU5d6KroDoSh4D6EzYGMGb%2BHa2QIYufRovLHuxJIL8FMzaqAeZ5p0nrBXmE0U%2FnRSAT7f3iyGV9txJwxJQY%2FXUnxkbByCoRC9rNlhf%2F9erKQ5qJ%2FrGLJPqaz6QGAB5POLu7wE0R4EmNoZcMggUYhWTxdEOCsQPOqy5mdWcwi4YDM%2BUNSJg0fCupifcRSn9GPXz5Rhz0myrQ9SnEq3deBdUfg92P4JE7o7XZbkVqieAxRcOHoSF%2FAOGg8T0OXY0F1s%2FPBAduiSHvP19FwBPiisuF5ZgBtHjqOyfh3Duraybtvxch5y%2Fkp7bTuqlqOPtFa1IzN410wkXEnb6XSY8eyZAaiaarhwAVgfsZp8jwDdTKecsr%2FF0i2md%2FGLUpgbi3U51ytVSpWxGwmDUkvUNplOujNszrE9N%2FX8aBk9m9RXeTBdBkT9%2FnkrxdGinOx%2FnsmHzzqo6chnLqn3DzRNGW5DrfITQkPevBm%2FJ0MPr6z6lEhHhIrBX31HjKnsyiefxjEUChr%2B4rP7eWYHXTVbUaWtWewypqus0X5SM8BjL%2BnPRo%2BX9N4RdAkWuLmQmfRN0yg10QJncWCISX998oUiGh3sUjCXt8vCBjqmARgK47NlLn6HXdxPL%2BXKzcDwZW2VICSX0vCquNz6fj1cQcuCW8d7bAGXRfYxQ%2Fvi49IBrks9pfjAKJePZXuq%2FGGMjhTO%2BnpcQCNIlfp7FWqYUZ9VhfCbECq7YabCidAHes9LKKuTJc7gdcb%2Bde8YFSpdn5CoDa7UPVK%2FOQ%2FdrsKSzstvLN71zPnDAx3aXFaTmuCrWWqEX6m7oKvTypuxG7uYDFITSIg%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2000%3A04%3A42%20GMT&X-Amz-Signature=9ebe4db533bb302690045c35b0924eeccfcce7c22a3f95fee1eb9ccabba17d68'
$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1750183612496/TAdN8NokD3-camera-video-segment-1750183614996.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNED233MW%2F20250618%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250618T160442Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIBxOo0S7hIeADNS8X%2F11v%2FFXQ5yY40M9wwoRF9Veb0znAiA0hEdIMmymyjsOlVstbTPCrvZ04Qw671YfXkdetNsKSCrmAwiR%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM14o7WDIZqelIU5d6KroDoSh4D6EzYGMGb%2BHa2QIYufRovLHuxJIL8FMzaqAeZ5p0nrBXmE0U%2FnRSAT7f3iyGV9txJwxJQY%2FXUnxkbByCoRC9rNlhf%2F9erKQ5qJ%2FrGLJPqaz6QGAB5POLu7wE0R4EmNoZcMggUYhWTxdEOCsQPOqy5mdWcwi4YDM%2BUNSJg0fCupifcRSn9GPXz5Rhz0myrQ9SnEq3deBdUfg92P4JE7o7XZbkVqieAxRcOHoSF%2FAOGg8T0OXY0F1s%2FPBAduiSHvP19FwBPiisuF5ZgBtHjqOyfh3Duraybtvxch5y%2Fkp7bTuqlqOPtFa1IzN410wkXEnb6XSY8eyZAaiaarhwAVgfsZp8jwDdTKecsr%2FF0i2md%2FGLUpgbi3U51ytVSpWxGwmDUkvUNplOujNszrE9N%2FX8aBk9m9RXeTBdBkT9%2FnkrxdGinOx%2FnsmHzzqo6chnLqn3DzRNGW5DrfITQkPevBm%2FJ0MPr6z6lEhHhIrBX31HjKnsyiefxjEUChr%2B4rP7eWYHXTVbUaWtWewypqus0X5SM8BjL%2BnPRo%2BX9N4RdAkWuLmQmfRN0yg10QJncWCISX998oUiGh3sUjCXt8vCBjqmARgK47NlLn6HXdxPL%2BXKzcDwZW2VICSX0vCquNz6fj1cQcuCW8d7bAGXRfYxQ%2Fvi49IBrks9pfjAKJePZXuq%2FGGMjhTO%2BnpcQCNIlfp7FWqYUZ9VhfCbECq7YabCidAHes9LKKuTJc7gdcb%2Bde8YFSpdn5CoDa7UPVK%2FOQ%2FdrsKSzstvLN71zPnDAx3aXFaTmuCrWWqEX6m7oKvTypuxG7uYDFITSIg%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2000%3A04%3A42%20GMT&X-Amz-Signature=4f1736ce16b63a5d5c67d7f60308d5b9f76ca7120acc29f1785f786b3674bae7'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991152736/1750182961683/iBjUYEdqdu-camera-video-segment-driver-1750182964183.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNED233MW%2F20250618%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250618T160442Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIBxOo0S7hIeADNS8X%2F11v%2FFXQ5yY40M9wwoRF9Veb0znAiA0hEdIMmymyjsOlVstbTPCrvZ04Qw671YfXkdetNsKSCrmAwiR%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM14o7WDIZqelIU5d6KroDoSh4D6EzYGMGb%2BHa2QIYufRovLHuxJIL8FMzaqAeZ5p0nrBXmE0U%2FnRSAT7f3iyGV9txJwxJQY%2FXUnxkbByCoRC9rNlhf%2F9erKQ5qJ%2FrGLJPqaz6QGAB5POLu7wE0R4EmNoZcMggUYhWTxdEOCsQPOqy5mdWcwi4YDM%2BUNSJg0fCupifcRSn9GPXz5Rhz0myrQ9SnEq3deBdUfg92P4JE7o7XZbkVqieAxRcOHoSF%2FAOGg8T0OXY0F1s%2FPBAduiSHvP19FwBPiisuF5ZgBtHjqOyfh3Duraybtvxch5y%2Fkp7bTuqlqOPtFa1IzN410wkXEnb6XSY8eyZAaiaarhwAVgfsZp8jwDdTKecsr%2FF0i2md%2FGLUpgbi3U51ytVSpWxGwmDUkvUNplOujNszrE9N%2FX8aBk9m9RXeTBdBkT9%2FnkrxdGinOx%2FnsmHzzqo6chnLqn3DzRNGW5DrfITQkPevBm%2FJ0MPr6z6lEhHhIrBX31HjKnsyiefxjEUChr%2B4rP7eWYHXTVbUaWtWewypqus0X5SM8BjL%2BnPRo%2BX9N4RdAkWuLmQmfRN0yg10QJncWCISX998oUiGh3sUjCXt8vCBjqmARgK47NlLn6HXdxPL%2BXKzcDwZW2VICSX0vCquNz6fj1cQcuCW8d7bAGXRfYxQ%2Fvi49IBrks9pfjAKJePZXuq%2FGGMjhTO%2BnpcQCNIlfp7FWqYUZ9VhfCbECq7YabCidAHes9LKKuTJc7gdcb%2Bde8YFSpdn5CoDa7UPVK%2FOQ%2FdrsKSzstvLN71zPnDAx3aXFaTmuCrWWqEX6m7oKvTypuxG7uYDFITSIg%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2000%3A04%3A42%20GMT&X-Amz-Signature=5cad44980d1a64cb4d66d6300ff7bc1ace6b9a003b0b49cd3ffcbebc14a7994c'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991152736/1750182961683/Foph7ZSROm-camera-video-segment-1750182964183.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNED233MW%2F20250618%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250618T160442Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIBxOo0S7hIeADNS8X%2F11v%2FFXQ5yY40M9wwoRF9Veb0znAiA0hEdIMmymyjsOlVstbTPCrvZ04Qw671YfXkdetNsKSCrmAwiR%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM14o7WDIZqelIU5d6KroDoSh4D6EzYGMGb%2BHa2QIYufRovLHuxJIL8FMzaqAeZ5p0nrBXmE0U%2FnRSAT7f3iyGV9txJwxJQY%2FXUnxkbByCoRC9rNlhf%2F9erKQ5qJ%2FrGLJPqaz6QGAB5POLu7wE0R4EmNoZcMggUYhWTxdEOCsQPOqy5mdWcwi4YDM%2BUNSJg0fCupifcRSn9GPXz5Rhz0myrQ9SnEq3deBdUfg92P4JE7o7XZbkVqieAxRcOHoSF%2FAOGg8T0OXY0F1s%2FPBAduiSHvP19FwBPiisuF5ZgBtHjqOyfh3Duraybtvxch5y%2Fkp7bTuqlqOPtFa1IzN410wkXEnb6XSY8eyZAaiaarhwAVgfsZp8jwDdTKecsr%2FF0i2md%2FGLUpgbi3U51ytVSpWxGwmDUkvUNplOujNszrE9N%2FX8aBk9m9RXeTBdBkT9%2FnkrxdGinOx%2FnsmHzzqo6chnLqn3DzRNGW5DrfITQkPevBm%2FJ0MPr6z6lEhHhIrBX31HjKnsyiefxjEUChr%2B4rP7eWYHXTVbUaWtWewypqus0X5SM8BjL%2BnPRo%2BX9N4RdAkWuLmQmfRN0yg10QJncWCISX998oUiGh3sUjCXt8vCBjqmARgK47NlLn6HXdxPL%2BXKzcDwZW2VICSX0vCquNz6fj1cQcuCW8d7bAGXRfYxQ%2Fvi49IBrks9pfjAKJePZXuq%2FGGMjhTO%2BnpcQCNIlfp7FWqYUZ9VhfCbECq7YabCidAHes9LKKuTJc7gdcb%2Bde8YFSpdn5CoDa7UPVK%2FOQ%2FdrsKSzstvLN71zPnDAx3aXFaTmuCrWWqEX6m7oKvTypuxG7uYDFITSIg%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2000%3A04%3A42%20GMT&X-Amz-Signature=20acb1b4525b331e9fb95684ffbbd620698c5a7e84b974b1ecfacb4daf5e9360'
